$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: predidx / pred_name update
$ws.Range("D20").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E20").Value = "['Normal']"

# Row 26: predidx / pred_name update
$ws.Range("D26").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E26").Value = "[]"

# Row 40: predidx / pred_name update
$ws.Range("D40").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E40").Value = "['Normal', 'HardwareFault']"
